$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "dct:modified" timestamp in B21
$ws.Range("B21").Value = "2023-08-17T12:13:29+00:00"

# Update skos:broader references in column F from "vocab.NNNN" to "vocab:NNNN"
# (dot separator replaced with colon separator to match new .ttl export)
$ws.Range("F24").Value = "vocab:1000"
$ws.Range("F25").Value = "vocab:1000"
$ws.Range("F26").Value = "vocab:1000"
$ws.Range("F27").Value = "vocab:1000,vocab:1007"
$ws.Range("F28").Value = "vocab:1000,vocab:1007"
$ws.Range("F29").Value = "vocab:1000,vocab:1007"
$ws.Range("F31").Value = "vocab:1007"
$ws.Range("F32").Value = "vocab:1007"
$ws.Range("F33").Value = "vocab:1007"
$ws.Range("F34").Value = "vocab:1007"
$ws.Range("F35").Value = "vocab:1007"
$ws.Range("F36").Value = "vocab:1007"
